$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.2444541603326797
$ws.Range("B1").Value = 0.2060756683349609
$ws.Range("C1").Value = 0.1889204978942871
$ws.Range("D1").Value = 0.2145034372806549
$ws.Range("E1").Value = 0.278813898563385
